# Apply updated crypto price/volume data (and the Kaspa/TheGraph row swap)
# to match the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'69.110.64"
$ws.Range("E2").Formula = "'  -3.79%  "
$ws.Range("D3").Formula = "'3.516.28"
$ws.Range("E3").Formula = "'  -4.61%  "
$ws.Range("D4").Formula = "'1.00"
$ws.Range("E4").Formula = "'  -0.05%  "
$ws.Range("D5").Formula = "'579.93"
$ws.Range("E5").Formula = "'  -1.45%  "
$ws.Range("D6").Formula = "'174.68"
$ws.Range("E6").Formula = "'  -2.36%  "
$ws.Range("D7").Formula = "'0.624"
$ws.Range("E7").Formula = "'  +0.15%  "
$ws.Range("D8").Formula = "'3.509.77"
$ws.Range("E8").Formula = "'  -4.54%  "
$ws.Range("E9").Formula = "'  +0.00%  "
$ws.Range("E10").Formula = "'  -6.68%  "
$ws.Range("E11").Formula = "'  +8.98%  "
$ws.Range("D12").Formula = "'0.602"
$ws.Range("E12").Formula = "'  -2.10%  "
$ws.Range("D13").Formula = "'47.34"
$ws.Range("E13").Formula = "'  -5.29%  "
$ws.Range("E14").Formula = "'  -3.25%  "
$ws.Range("D15").Formula = "'671.80"
$ws.Range("E15").Formula = "'  -1.72%  "
$ws.Range("D16").Formula = "'4.080.99"
$ws.Range("E16").Formula = "'  -4.76%  "
$ws.Range("D17").Formula = "'8.85"
$ws.Range("E17").Formula = "'  -1.25%  "
$ws.Range("D18").Formula = "'3.510.21"
$ws.Range("E18").Formula = "'  -5.02%  "
$ws.Range("D19").Formula = "'68.970.72"
$ws.Range("E19").Formula = "'  -4.05%  "
$ws.Range("E20").Formula = "'  -1.58%  "
$ws.Range("E21").Formula = "'  -3.45%  "
$ws.Range("D22").Formula = "'11.22"
$ws.Range("E22").Formula = "'  -3.59%  "
$ws.Range("D23").Formula = "'0.907"
$ws.Range("E23").Formula = "'  -3.55%  "
$ws.Range("D24").Formula = "'16.31"
$ws.Range("E24").Formula = "'  -8.26%  "
$ws.Range("D25").Formula = "'98.47"
$ws.Range("E25").Formula = "'  -4.89%  "
$ws.Range("D26").Formula = "'3.88"
$ws.Range("E26").Formula = "'  -4.10%  "
$ws.Range("E27").Formula = "'  -0.79%  "
$ws.Range("E28").Formula = "'  +0.14%  "
$ws.Range("D29").Formula = "'2.66"
$ws.Range("E29").Formula = "'  -6.77%  "
$ws.Range("D30").Formula = "'9.45"
$ws.Range("E30").Formula = "'  -7.32%  "
$ws.Range("D31").Formula = "'32.99"
$ws.Range("E31").Formula = "'  -7.02%  "
$ws.Range("D32").Formula = "'3.23"
$ws.Range("E32").Formula = "'  -7.32%  "
$ws.Range("D33").Formula = "'8.77"
$ws.Range("E33").Formula = "'  -4.73%  "
$ws.Range("D34").Formula = "'7.34"
$ws.Range("E34").Formula = "'  -0.49%  "
$ws.Range("E35").Formula = "'  -4.84%  "
$ws.Range("D36").Formula = "'577.85"
$ws.Range("E36").Formula = "'  +0.60%  "
$ws.Range("D37").Formula = "'10.97"
$ws.Range("E37").Formula = "'  -3.06%  "
$ws.Range("D38").Formula = "'3.60"
$ws.Range("E38").Formula = "'  -14.34%  "
$ws.Range("E39").Formula = "'  -3.75%  "
$ws.Range("D40").Formula = "'57.11"
$ws.Range("E40").Formula = "'  -4.44%  "
$ws.Range("D41").Formula = "'0.998"
$ws.Range("E41").Formula = "'  +0.01%  "
$ws.Range("B42").Formula = "TheGraph"
$ws.Range("C42").Formula = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Formula = "'0.338"
$ws.Range("E42").Formula = "'  -2.98%  "
$ws.Range("B43").Formula = "Kaspa"
$ws.Range("C43").Formula = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Formula = "'0.137"
$ws.Range("E43").Formula = "'  -5.39%  "
$ws.Range("D44").Formula = "'0.0439"
$ws.Range("E44").Formula = "'  -5.31%  "
$ws.Range("D45").Formula = "'3.427.03"
$ws.Range("E45").Formula = "'  -8.54%  "
$ws.Range("D46").Formula = "'33.47"
$ws.Range("E46").Formula = "'  -5.66%  "
$ws.Range("D47").Formula = "'0.0₃0708"
$ws.Range("E47").Formula = "'  -8.41%  "
$ws.Range("D48").Formula = "'2.93"
$ws.Range("E48").Formula = "'  +1.86%  "
$ws.Range("D49").Formula = "'2.60"
$ws.Range("E49").Formula = "'  -6.90%  "
$ws.Range("E50").Formula = "'  -0.56%  "
$ws.Range("D51").Formula = "'131.05"
$ws.Range("E51").Formula = "'  -2.24%  "
